$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2026-01-14T15:34:52+00:00"
$meta.Range("B12").Value = "Acte"

# --- Elements sheet ---
$elem = $wb.Worksheets.Item("Elements")
$elem.Range("M2").Value = "Acte"

$elem.Range("L13").Value = "Auteur"
$elem.Range("M13").Value = "Auteur"

$elem.Range("L14").Value = "Informateur"
$elem.Range("M14").Value = "Informateur"

$elem.Range("L15").Value = "Participant"
$elem.Range("M15").Value = "Participant"
